# Apply the GroPIN metadataschema23 update: add author/publication metadata
# for rows 4-7, 15, 28, 34, 35, 39, and rework the secondary-model parameter
# rows 133-137 (new real parameter descriptions instead of placeholder
# "dummy"/"no idea" values, and dropping the now-unused Source/Subject/Dist
# columns for the Constant rows, collapsing rows 136-137 down to a single
# "Add" hint in column P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: corresponding author (Marcel Fuhrmann) + first paper author (Devlieghere, F.) ---
$ws.Range("M4").Value = "Marcel"
$ws.Range("O4").Value = "Fuhrmann"
$ws.Range("P4").Value = "BfR"
$ws.Range("R4").Value = "marcel.fuhrmann@bfr.bund.de"
$ws.Range("AC4").Value = "F."
$ws.Range("AE4").Value = "Devlieghere"
$ws.Range("AH4").Value = "none given"

# --- Row 5: paper author (Lefevere, I.) ---
$ws.Range("AC5").Value = "I."
$ws.Range("AE5").Value = "Lefevere"
$ws.Range("AH5").Value = "none given"

# --- Row 6: paper author (Magnin, A.) ---
$ws.Range("AC6").Value = "A."
$ws.Range("AE6").Value = "Magnin"
$ws.Range("AH6").Value = "none given"

# --- Row 7: paper author (Debevere, J.) ---
$ws.Range("AC7").Value = "J."
$ws.Range("AE7").Value = "Debevere"
$ws.Range("AH7").Value = "none given"

# --- Row 15: publication reference ---
$ws.Range("K15").Value = "Yes"
$ws.Range("O15").Value = "none given"
$ws.Range("P15").Value = "Devlieghere,  F.,  Lefevere,  I.,  Magnin,  A.,  Debevere,  J."
$ws.Range("Q15").Value = "Growth of Aeromonas hydrophila in modified-atmosphere-packed cooked meat products"
$ws.Range("S15").Value = "Food Microbiology ,  17,185-196"

# --- Row 28: model category ---
$ws.Range("I28").Value = "QRA model"

# --- Row 34: objective ---
$ws.Range("I34").Value = "This model predicts and visualize the mu_max of Aeromonas hydrophila in modified BHI with the independent variable(s) T, aw according to the publication from Devlieghere, F., Lefevere, I., Magnin, A., Debevere, J., 2000 on Growth of Aeromonas hydrophila in modified-atmosphere-packed cooked meat products. "

# --- Row 35: description ---
$ws.Range("I35").Value = "This model and all metadata included have been automatically generated from the GroPIN microbial modelling DataBase (https://www.aua.gr/psomas/gropin/, version 2020). The model code has been converted from Excel to R and the model itself is provided as an FSKX file. This FSKX model contains also an R script to visualize model-based prediction results similar to those visualizations provided by the GroPIN software. A user of the FSKX model can provide user-defined values for all model input parameters, some of them specifically introduced to customize the generated visualization."

# --- Row 39: scope / matrix ---
$ws.Range("K39").Value = "modified BHI"
$ws.Range("M39").Value = "none given"
$ws.Range("W39").Value = "Aeromonas hydrophila"
$ws.Range("Y39").Value = "log10(CFU)"

# --- Row 133: input parameter "T" -> real Temperature metadata ---
$ws.Range("N133").Value = "Temperature"
$ws.Range("O133").Value = "descr Temperature"
$ws.Range("P133").Value = "C"
$ws.Range("Q133").Value = "Temperature"
$ws.Range("S133:U133").ClearContents()
$ws.Range("V133").Value = "seq(3.996003996004,12.012,length.out=21)"

# --- Row 134: input parameter "aw" -> real water activity metadata ---
$ws.Range("N134").Value = "water activity"
$ws.Range("O134").Value = "descr water activity"
$ws.Range("P134").Value = "[%]"
$ws.Range("Q134").Value = "Dimensionless Parameter"
$ws.Range("S134:U134").ClearContents()
$ws.Range("V134").Value = "seq(0.973026973026973,0.992992,length.out=21)"

# --- Row 135: constant "a" -> output "mumax" data frame ---
$ws.Range("L135").Value = "mumax"
$ws.Range("M135").Value = "Output"
$ws.Range("N135").Value = "data frame with variables and corresponding mumax"
$ws.Range("O135").Value = "This dataframe consists of a number of columns " + [char]10 + "                            in relation to the number of variables of this" + [char]10 + "                            model. One additional column contains the response" + [char]10 + "                            surface mu_max result based on this secondary model."
$ws.Range("P135").Value = "[]"
$ws.Range("Q135").ClearContents()
$ws.Range("R135").Value = "Matrix[number,number]"
$ws.Range("S135:V135").ClearContents()

# --- Row 136: constant "awmin" entirely dropped, only a "Add" hint remains ---
$ws.Range("L136:O136").ClearContents()
$ws.Range("Q136:V136").ClearContents()
$ws.Range("P136").Value = "Add"

# --- Row 137: constant "Tmin" entirely dropped, only a "Add" hint remains ---
$ws.Range("L137:O137").ClearContents()
$ws.Range("Q137:V137").ClearContents()
$ws.Range("P137").Value = "Add"
